# Re-label the J4 "A" bus pins to their correct (new) silkscreen pin names,
# re-order / relabel the J4 "B" bus pins, fix the J3 clock pins to upper-case
# silkscreen names, move the CLK_EXT_IN / CLK_OUT rows from J3 to J4, and add
# the two new SW0 / SW1 rows documenting the switch inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J4 "A" bus (rows 6-15): new pin numbers + new silkscreen pin names ---
$ws.Range("C6").Value  = 24
$ws.Range("D6").Value  = "H15"

$ws.Range("C7").Value  = 20
$ws.Range("D7").Value  = "J15"

$ws.Range("C8").Value  = 17
$ws.Range("D8").Value  = "K15"

$ws.Range("C9").Value  = 16
$ws.Range("D9").Value  = "K16"

$ws.Range("C10").Value = 15
$ws.Range("D10").Value = "K14"

$ws.Range("C11").Value = 14
$ws.Range("D11").Value = "L16"

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = "M16"

$ws.Range("C13").Value = 8
$ws.Range("D13").Value = "N16"

$ws.Range("C14").Value = 7
$ws.Range("D14").Value = "L14"

$ws.Range("C15").Value = 6
$ws.Range("D15").Value = "N14"

# --- J4 "B" bus (rows 17-26): new pin numbers + new silkscreen pin names ---
$ws.Range("C17").Value = 39
$ws.Range("D17").Value = "G14"

$ws.Range("C18").Value = 37
$ws.Range("D18").Value = "D14"

$ws.Range("C19").Value = 38
$ws.Range("D19").Value = "F14"

$ws.Range("C20").Value = 35
$ws.Range("D20").Value = "D16"

$ws.Range("C21").Value = 34
$ws.Range("D21").Value = "E14"

$ws.Range("C22").Value = 33
$ws.Range("D22").Value = "E16"

$ws.Range("C23").Value = 30
$ws.Range("D23").Value = "E15"

$ws.Range("C24").Value = 29
$ws.Range("D24").Value = "F16"

$ws.Range("C25").Value = 28
$ws.Range("D25").Value = "F15"

$ws.Range("C26").Value = 26
$ws.Range("D26").Value = "G15"

# --- J3 clock pins (rows 28-33): pin numbers unchanged, pin names upper-cased ---
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = "D10"

$ws.Range("D29").Value = "E10"
$ws.Range("D30").Value = "B9"
$ws.Range("D31").Value = "A10"
$ws.Range("D32").Value = "D8"
$ws.Range("D33").Value = "E9"

# --- ~CW moves from J4 pin 9 / m15 to J4 pin 5 / M14 ---
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = "M14"

# --- CLK_EXT_IN / CLK_OUT move from connector J3 to connector J4, pin names upper-cased ---
$ws.Range("B37").Value = "J4"
$ws.Range("D37").Value = "H16"

$ws.Range("B38").Value = "J4"
$ws.Range("D38").Value = "H14"

# --- New rows documenting the switch inputs on J3 ---
$ws.Range("B45").Value = "J3"
$ws.Range("C45").Value = 5
$ws.Range("E45").Value = "SW0"

$ws.Range("B46").Value = "J3"
$ws.Range("C46").Value = 6
$ws.Range("E46").Value = "SW1"

$ws.Range("D45").Value = "F8"
$ws.Range("D46").Value = "B12"

# --- Refresh the view: scroll up a bit and move the selection (cosmetic) ---
$ws.Range("G37").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
